$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 10011
$ws.Cells.Item(2, 4).Value = 13744166
$ws.Cells.Item(4, 3).Value = 19780
$ws.Cells.Item(4, 4).Value = 25345251
$ws.Cells.Item(6, 3).Value = 1083
$ws.Cells.Item(6, 4).Value = 1346437
$ws.Cells.Item(7, 3).Value = 55763
$ws.Cells.Item(7, 4).Value = 80224387
$ws.Cells.Item(8, 3).Value = 71443
$ws.Cells.Item(8, 4).Value = 95277876
$ws.Cells.Item(9, 3).Value = 19327
$ws.Cells.Item(9, 4).Value = 26879754
$ws.Cells.Item(10, 3).Value = 57595
$ws.Cells.Item(10, 4).Value = 82660917
$ws.Cells.Item(11, 3).Value = 7714
$ws.Cells.Item(11, 4).Value = 10097573
$ws.Cells.Item(12, 3).Value = 3470
$ws.Cells.Item(12, 4).Value = 4893511
$ws.Cells.Item(13, 3).Value = 13103
$ws.Cells.Item(13, 4).Value = 17770692
$ws.Cells.Item(14, 3).Value = 40972
$ws.Cells.Item(14, 4).Value = 55412097
$ws.Cells.Item(15, 3).Value = 19226
$ws.Cells.Item(15, 4).Value = 24911486
$ws.Cells.Item(17, 3).Value = 38179
$ws.Cells.Item(17, 4).Value = 48618451
$ws.Cells.Item(18, 3).Value = 45730
$ws.Cells.Item(18, 4).Value = 61323887
$ws.Cells.Item(19, 3).Value = 22730
$ws.Cells.Item(19, 4).Value = 27831985
$ws.Cells.Item(20, 3).Value = 47034
$ws.Cells.Item(20, 4).Value = 57429044
$ws.Cells.Item(21, 3).Value = 3426
$ws.Cells.Item(21, 4).Value = 4714925
$ws.Cells.Item(23, 3).Value = 5521
$ws.Cells.Item(23, 4).Value = 6961638
$ws.Cells.Item(25, 3).Value = 446
$ws.Cells.Item(25, 4).Value = 557774
$ws.Cells.Item(26, 3).Value = 13759
$ws.Cells.Item(26, 4).Value = 19722073
$ws.Cells.Item(27, 3).Value = 22104
$ws.Cells.Item(27, 4).Value = 29072778
$ws.Cells.Item(28, 3).Value = 2660
$ws.Cells.Item(28, 4).Value = 3617158
$ws.Cells.Item(29, 3).Value = 16721
$ws.Cells.Item(29, 4).Value = 23831237
$ws.Cells.Item(30, 3).Value = 1424
$ws.Cells.Item(30, 4).Value = 1780528
$ws.Cells.Item(31, 3).Value = 826
$ws.Cells.Item(31, 4).Value = 1127208
$ws.Cells.Item(32, 3).Value = 2931
$ws.Cells.Item(32, 4).Value = 3878264
$ws.Cells.Item(33, 3).Value = 7964
$ws.Cells.Item(33, 4).Value = 10747558
$ws.Cells.Item(34, 3).Value = 3904
$ws.Cells.Item(34, 4).Value = 4970430
$ws.Cells.Item(36, 3).Value = 5647
$ws.Cells.Item(36, 4).Value = 6785304
$ws.Cells.Item(37, 3).Value = 10044
$ws.Cells.Item(37, 4).Value = 13373148
$ws.Cells.Item(38, 3).Value = 5744
$ws.Cells.Item(38, 4).Value = 6861910
$ws.Cells.Item(39, 3).Value = 14389
$ws.Cells.Item(39, 4).Value = 17741218
$ws.Cells.Item(40, 3).Value = 2820
$ws.Cells.Item(40, 4).Value = 3857678
$ws.Cells.Item(42, 3).Value = 7343
$ws.Cells.Item(42, 4).Value = 9047696
$ws.Cells.Item(44, 3).Value = 292
$ws.Cells.Item(44, 4).Value = 348983
$ws.Cells.Item(45, 3).Value = 15381
$ws.Cells.Item(45, 4).Value = 21985163
$ws.Cells.Item(46, 3).Value = 23924
$ws.Cells.Item(46, 4).Value = 32069662
$ws.Cells.Item(47, 3).Value = 3318
$ws.Cells.Item(47, 4).Value = 4576611
$ws.Cells.Item(48, 3).Value = 22764
$ws.Cells.Item(48, 4).Value = 32655928
$ws.Cells.Item(49, 3).Value = 2303
$ws.Cells.Item(49, 4).Value = 2905775
$ws.Cells.Item(50, 3).Value = 1047
$ws.Cells.Item(50, 4).Value = 1475460
$ws.Cells.Item(51, 3).Value = 4340
$ws.Cells.Item(51, 4).Value = 5665590
$ws.Cells.Item(52, 3).Value = 11830
$ws.Cells.Item(52, 4).Value = 15913727
$ws.Cells.Item(53, 3).Value = 5011
$ws.Cells.Item(53, 4).Value = 6225148
$ws.Cells.Item(55, 3).Value = 6338
$ws.Cells.Item(55, 4).Value = 7767297
$ws.Cells.Item(56, 3).Value = 16779
$ws.Cells.Item(56, 4).Value = 22646697
$ws.Cells.Item(57, 3).Value = 7178
$ws.Cells.Item(57, 4).Value = 8591159
$ws.Cells.Item(58, 3).Value = 16103
$ws.Cells.Item(58, 4).Value = 19941989
$ws.Cells.Item(59, 3).Value = 2504
$ws.Cells.Item(59, 4).Value = 3437542
$ws.Cells.Item(60, 3).Value = 4625
$ws.Cells.Item(60, 4).Value = 5888006
$ws.Cells.Item(63, 3).Value = 13551
$ws.Cells.Item(63, 4).Value = 19365772
$ws.Cells.Item(64, 3).Value = 17918
$ws.Cells.Item(64, 4).Value = 23396186
$ws.Cells.Item(65, 3).Value = 3604
$ws.Cells.Item(65, 4).Value = 5071157
$ws.Cells.Item(66, 3).Value = 13339
$ws.Cells.Item(66, 4).Value = 19129935
$ws.Cells.Item(67, 3).Value = 1552
$ws.Cells.Item(67, 4).Value = 2015412
$ws.Cells.Item(68, 3).Value = 699
$ws.Cells.Item(68, 4).Value = 975899
$ws.Cells.Item(69, 3).Value = 2970
$ws.Cells.Item(69, 4).Value = 3986902
$ws.Cells.Item(70, 3).Value = 7436
$ws.Cells.Item(70, 4).Value = 10060824
$ws.Cells.Item(71, 3).Value = 3990
$ws.Cells.Item(71, 4).Value = 5047018
$ws.Cells.Item(73, 3).Value = 4588
$ws.Cells.Item(73, 4).Value = 5713908
$ws.Cells.Item(74, 3).Value = 8570
$ws.Cells.Item(74, 4).Value = 11349452
$ws.Cells.Item(75, 3).Value = 5152
$ws.Cells.Item(75, 4).Value = 6303521
$ws.Cells.Item(76, 3).Value = 12905
$ws.Cells.Item(76, 4).Value = 15889917
$ws.Cells.Item(77, 3).Value = 2275
$ws.Cells.Item(77, 4).Value = 3135984
$ws.Cells.Item(78, 3).Value = 1789
$ws.Cells.Item(78, 4).Value = 2366386
$ws.Cells.Item(79, 3).Value = 29
$ws.Cells.Item(79, 4).Value = 42026
$ws.Cells.Item(80, 3).Value = 4478
$ws.Cells.Item(80, 4).Value = 6357876
$ws.Cells.Item(81, 3).Value = 4750
$ws.Cells.Item(81, 4).Value = 6613105
$ws.Cells.Item(82, 3).Value = 739
$ws.Cells.Item(82, 4).Value = 1052507
$ws.Cells.Item(83, 3).Value = 4885
$ws.Cells.Item(83, 4).Value = 7015590
$ws.Cells.Item(84, 3).Value = 316
$ws.Cells.Item(84, 4).Value = 422834
$ws.Cells.Item(86, 3).Value = 1180
$ws.Cells.Item(86, 4).Value = 1610888
$ws.Cells.Item(87, 3).Value = 3077
$ws.Cells.Item(87, 4).Value = 4290367
$ws.Cells.Item(88, 3).Value = 1734
$ws.Cells.Item(88, 4).Value = 2210537
$ws.Cells.Item(89, 3).Value = 882
$ws.Cells.Item(89, 4).Value = 1080606
$ws.Cells.Item(90, 3).Value = 1711
$ws.Cells.Item(90, 4).Value = 2301028
$ws.Cells.Item(91, 3).Value = 984
$ws.Cells.Item(91, 4).Value = 1281014
$ws.Cells.Item(92, 3).Value = 2570
$ws.Cells.Item(92, 4).Value = 3088118
$ws.Cells.Item(93, 3).Value = 4815
$ws.Cells.Item(93, 4).Value = 6713426
$ws.Cells.Item(95, 3).Value = 9664
$ws.Cells.Item(95, 4).Value = 12602123
$ws.Cells.Item(97, 3).Value = 1020
$ws.Cells.Item(97, 4).Value = 1246473
$ws.Cells.Item(98, 3).Value = 28478
$ws.Cells.Item(98, 4).Value = 40933184
$ws.Cells.Item(99, 3).Value = 40458
$ws.Cells.Item(99, 4).Value = 53782925
$ws.Cells.Item(100, 3).Value = 6875
$ws.Cells.Item(100, 4).Value = 9357627
$ws.Cells.Item(101, 3).Value = 27374
$ws.Cells.Item(101, 4).Value = 39464408
$ws.Cells.Item(102, 3).Value = 3145
$ws.Cells.Item(102, 4).Value = 4042757
$ws.Cells.Item(103, 3).Value = 1743
$ws.Cells.Item(103, 4).Value = 2420550
$ws.Cells.Item(104, 3).Value = 5075
$ws.Cells.Item(104, 4).Value = 6840326
$ws.Cells.Item(105, 3).Value = 17479
$ws.Cells.Item(105, 4).Value = 23498512
$ws.Cells.Item(106, 3).Value = 7647
$ws.Cells.Item(106, 4).Value = 9772432
$ws.Cells.Item(108, 3).Value = 9604
$ws.Cells.Item(108, 4).Value = 11719830
$ws.Cells.Item(109, 3).Value = 20607
$ws.Cells.Item(109, 4).Value = 28107719
$ws.Cells.Item(110, 3).Value = 9547
$ws.Cells.Item(110, 4).Value = 11381333
$ws.Cells.Item(111, 3).Value = 29321
$ws.Cells.Item(111, 4).Value = 35284260
$ws.Cells.Item(113, 3).Value = 4762
$ws.Cells.Item(113, 4).Value = 6593681
$ws.Cells.Item(114, 3).Value = 2532
$ws.Cells.Item(114, 4).Value = 3538719
$ws.Cells.Item(117, 3).Value = 4507
$ws.Cells.Item(117, 4).Value = 6512067
$ws.Cells.Item(118, 3).Value = 6863
$ws.Cells.Item(118, 4).Value = 9513872
$ws.Cells.Item(119, 3).Value = 1586
$ws.Cells.Item(119, 4).Value = 2232548
$ws.Cells.Item(120, 3).Value = 4480
$ws.Cells.Item(120, 4).Value = 6416176
$ws.Cells.Item(121, 3).Value = 607
$ws.Cells.Item(121, 4).Value = 836016
$ws.Cells.Item(123, 3).Value = 925
$ws.Cells.Item(123, 4).Value = 1288411
$ws.Cells.Item(124, 3).Value = 2483
$ws.Cells.Item(124, 4).Value = 3483454
$ws.Cells.Item(125, 3).Value = 2874
$ws.Cells.Item(125, 4).Value = 3928130
$ws.Cells.Item(126, 3).Value = 1448
$ws.Cells.Item(126, 4).Value = 1855173
$ws.Cells.Item(127, 3).Value = 2083
$ws.Cells.Item(127, 4).Value = 2943043
$ws.Cells.Item(128, 3).Value = 961
$ws.Cells.Item(128, 4).Value = 1291630
$ws.Cells.Item(129, 3).Value = 2522
$ws.Cells.Item(129, 4).Value = 3294320
$ws.Cells.Item(130, 3).Value = 1353
$ws.Cells.Item(130, 4).Value = 1840207
$ws.Cells.Item(132, 3).Value = 621
$ws.Cells.Item(132, 4).Value = 876276
$ws.Cells.Item(134, 3).Value = 1541
$ws.Cells.Item(134, 4).Value = 2254935
$ws.Cells.Item(135, 3).Value = 1517
$ws.Cells.Item(135, 4).Value = 2136048
$ws.Cells.Item(136, 3).Value = 480
$ws.Cells.Item(136, 4).Value = 702384
$ws.Cells.Item(137, 3).Value = 1025
$ws.Cells.Item(137, 4).Value = 1464422
$ws.Cells.Item(138, 3).Value = 127
$ws.Cells.Item(138, 4).Value = 180251
$ws.Cells.Item(141, 3).Value = 667
$ws.Cells.Item(141, 4).Value = 945642
$ws.Cells.Item(142, 3).Value = 671
$ws.Cells.Item(142, 4).Value = 952702
$ws.Cells.Item(143, 3).Value = 306
$ws.Cells.Item(143, 4).Value = 420111
$ws.Cells.Item(145, 3).Value = 240
$ws.Cells.Item(145, 4).Value = 331121
$ws.Cells.Item(146, 3).Value = 527
$ws.Cells.Item(146, 4).Value = 717118
$ws.Cells.Item(147, 3).Value = 2623
$ws.Cells.Item(147, 4).Value = 3576144
$ws.Cells.Item(149, 3).Value = 7351
$ws.Cells.Item(149, 4).Value = 9138828
$ws.Cells.Item(152, 3).Value = 23605
$ws.Cells.Item(152, 4).Value = 33526815
$ws.Cells.Item(153, 3).Value = 39014
$ws.Cells.Item(153, 4).Value = 49763066
$ws.Cells.Item(154, 3).Value = 10194
$ws.Cells.Item(154, 4).Value = 13977378
$ws.Cells.Item(155, 3).Value = 26408
$ws.Cells.Item(155, 4).Value = 38118173
$ws.Cells.Item(156, 3).Value = 3147
$ws.Cells.Item(156, 4).Value = 4072801
$ws.Cells.Item(157, 3).Value = 1774
$ws.Cells.Item(157, 4).Value = 2482241
$ws.Cells.Item(158, 3).Value = 4722
$ws.Cells.Item(158, 4).Value = 6340234
$ws.Cells.Item(159, 3).Value = 16853
$ws.Cells.Item(159, 4).Value = 22930790
$ws.Cells.Item(160, 3).Value = 7564
$ws.Cells.Item(160, 4).Value = 9454036
$ws.Cells.Item(162, 3).Value = 8556
$ws.Cells.Item(162, 4).Value = 10630479
$ws.Cells.Item(163, 3).Value = 21181
$ws.Cells.Item(163, 4).Value = 28774009
$ws.Cells.Item(164, 3).Value = 9426
$ws.Cells.Item(164, 4).Value = 11400894
$ws.Cells.Item(165, 3).Value = 27915
$ws.Cells.Item(165, 4).Value = 32818687
$ws.Cells.Item(166, 3).Value = 791
$ws.Cells.Item(166, 4).Value = 1092716
$ws.Cells.Item(168, 3).Value = 19564
$ws.Cells.Item(168, 4).Value = 26191334
$ws.Cells.Item(169, 3).Value = 222
$ws.Cells.Item(169, 4).Value = 318594
$ws.Cells.Item(170, 3).Value = 798
$ws.Cells.Item(170, 4).Value = 1088979
$ws.Cells.Item(171, 3).Value = 60149
$ws.Cells.Item(171, 4).Value = 86759768
$ws.Cells.Item(172, 3).Value = 104134
$ws.Cells.Item(172, 4).Value = 142637951
$ws.Cells.Item(173, 3).Value = 125718
$ws.Cells.Item(173, 4).Value = 179709158
$ws.Cells.Item(174, 3).Value = 66325
$ws.Cells.Item(174, 4).Value = 97559188
$ws.Cells.Item(175, 3).Value = 30791
$ws.Cells.Item(175, 4).Value = 42194852
$ws.Cells.Item(176, 3).Value = 7262
$ws.Cells.Item(176, 4).Value = 10399268
$ws.Cells.Item(177, 3).Value = 18638
$ws.Cells.Item(177, 4).Value = 26388123
$ws.Cells.Item(178, 3).Value = 117004
$ws.Cells.Item(178, 4).Value = 160953355
$ws.Cells.Item(179, 3).Value = 30702
$ws.Cells.Item(179, 4).Value = 41201789
$ws.Cells.Item(181, 3).Value = 32157
$ws.Cells.Item(181, 4).Value = 39730922
$ws.Cells.Item(182, 3).Value = 55433
$ws.Cells.Item(182, 4).Value = 75112352
$ws.Cells.Item(183, 3).Value = 47587
$ws.Cells.Item(183, 4).Value = 61151123
$ws.Cells.Item(184, 3).Value = 56949
$ws.Cells.Item(184, 4).Value = 74280936
$ws.Cells.Item(185, 3).Value = 3897
$ws.Cells.Item(185, 4).Value = 5158726
$ws.Cells.Item(187, 3).Value = 3587
$ws.Cells.Item(187, 4).Value = 4827351
$ws.Cells.Item(188, 3).Value = 10
$ws.Cells.Item(188, 4).Value = 14213
$ws.Cells.Item(190, 3).Value = 6859
$ws.Cells.Item(190, 4).Value = 9980326
$ws.Cells.Item(191, 3).Value = 11674
$ws.Cells.Item(191, 4).Value = 15963941
$ws.Cells.Item(192, 3).Value = 1614
$ws.Cells.Item(192, 4).Value = 2266755
$ws.Cells.Item(193, 3).Value = 5961
$ws.Cells.Item(193, 4).Value = 8451025
$ws.Cells.Item(194, 3).Value = 804
$ws.Cells.Item(194, 4).Value = 1085425
$ws.Cells.Item(196, 3).Value = 1332
$ws.Cells.Item(196, 4).Value = 1861837
$ws.Cells.Item(197, 3).Value = 3728
$ws.Cells.Item(197, 4).Value = 5241932
$ws.Cells.Item(198, 3).Value = 2135
$ws.Cells.Item(198, 4).Value = 2946799
$ws.Cells.Item(199, 3).Value = 2996
$ws.Cells.Item(199, 4).Value = 3917393
$ws.Cells.Item(200, 3).Value = 4896
$ws.Cells.Item(200, 4).Value = 6925603
$ws.Cells.Item(201, 3).Value = 1925
$ws.Cells.Item(201, 4).Value = 2525204
$ws.Cells.Item(202, 3).Value = 4595
$ws.Cells.Item(202, 4).Value = 5864317
$ws.Cells.Item(203, 3).Value = 1461
$ws.Cells.Item(203, 4).Value = 1840401
$ws.Cells.Item(204, 3).Value = 1709
$ws.Cells.Item(204, 4).Value = 2310341
$ws.Cells.Item(207, 3).Value = 2745
$ws.Cells.Item(207, 4).Value = 3970381
$ws.Cells.Item(208, 3).Value = 4562
$ws.Cells.Item(208, 4).Value = 6275303
$ws.Cells.Item(209, 3).Value = 1496
$ws.Cells.Item(209, 4).Value = 2116463
$ws.Cells.Item(210, 3).Value = 2635
$ws.Cells.Item(210, 4).Value = 3783913
$ws.Cells.Item(211, 3).Value = 428
$ws.Cells.Item(211, 4).Value = 575403
$ws.Cells.Item(212, 3).Value = 194
$ws.Cells.Item(212, 4).Value = 277378
$ws.Cells.Item(213, 3).Value = 562
$ws.Cells.Item(213, 4).Value = 788591
$ws.Cells.Item(214, 3).Value = 2109
$ws.Cells.Item(214, 4).Value = 2934961
$ws.Cells.Item(215, 3).Value = 1979
$ws.Cells.Item(215, 4).Value = 2700069
$ws.Cells.Item(216, 3).Value = 1084
$ws.Cells.Item(216, 4).Value = 1408335
$ws.Cells.Item(217, 3).Value = 1773
$ws.Cells.Item(217, 4).Value = 2485149
$ws.Cells.Item(218, 3).Value = 790
$ws.Cells.Item(218, 4).Value = 1060802
$ws.Cells.Item(219, 3).Value = 2332
$ws.Cells.Item(219, 4).Value = 3053430
$ws.Cells.Item(221, 3).Value = 1573
$ws.Cells.Item(221, 4).Value = 2315588
$ws.Cells.Item(222, 3).Value = 410
$ws.Cells.Item(222, 4).Value = 597653
$ws.Cells.Item(224, 3).Value = 1265
$ws.Cells.Item(224, 4).Value = 1870288
$ws.Cells.Item(225, 3).Value = 5550
$ws.Cells.Item(225, 4).Value = 7920269
$ws.Cells.Item(226, 3).Value = 1236
$ws.Cells.Item(226, 4).Value = 1822235
$ws.Cells.Item(227, 3).Value = 676
$ws.Cells.Item(227, 4).Value = 989454
$ws.Cells.Item(232, 3).Value = 239
$ws.Cells.Item(232, 4).Value = 349464
$ws.Cells.Item(237, 3).Value = 2922
$ws.Cells.Item(237, 4).Value = 3976444
$ws.Cells.Item(239, 3).Value = 5142
$ws.Cells.Item(239, 4).Value = 6541964
$ws.Cells.Item(242, 3).Value = 14634
$ws.Cells.Item(242, 4).Value = 20946386
$ws.Cells.Item(243, 3).Value = 24691
$ws.Cells.Item(243, 4).Value = 32364180
$ws.Cells.Item(244, 3).Value = 3765
$ws.Cells.Item(244, 4).Value = 5183541
$ws.Cells.Item(245, 3).Value = 18192
$ws.Cells.Item(245, 4).Value = 26140226
$ws.Cells.Item(246, 3).Value = 1686
$ws.Cells.Item(246, 4).Value = 2132650
$ws.Cells.Item(247, 3).Value = 1089
$ws.Cells.Item(247, 4).Value = 1511397
$ws.Cells.Item(248, 3).Value = 3438
$ws.Cells.Item(248, 4).Value = 4547581
$ws.Cells.Item(249, 3).Value = 9920
$ws.Cells.Item(249, 4).Value = 13437567
$ws.Cells.Item(250, 3).Value = 4473
$ws.Cells.Item(250, 4).Value = 5611883
$ws.Cells.Item(252, 3).Value = 5484
$ws.Cells.Item(252, 4).Value = 6608455
$ws.Cells.Item(253, 3).Value = 10174
$ws.Cells.Item(253, 4).Value = 13474783
$ws.Cells.Item(254, 3).Value = 6238
$ws.Cells.Item(254, 4).Value = 7619694
$ws.Cells.Item(255, 3).Value = 17036
$ws.Cells.Item(255, 4).Value = 20759644
$ws.Cells.Item(256, 3).Value = 10838
$ws.Cells.Item(256, 4).Value = 14986495
$ws.Cells.Item(258, 3).Value = 15775
$ws.Cells.Item(258, 4).Value = 19486991
$ws.Cells.Item(260, 3).Value = 816
$ws.Cells.Item(260, 4).Value = 924020
$ws.Cells.Item(261, 3).Value = 42171
$ws.Cells.Item(261, 4).Value = 59882077
$ws.Cells.Item(262, 3).Value = 59335
$ws.Cells.Item(262, 4).Value = 78101094
$ws.Cells.Item(263, 3).Value = 9433
$ws.Cells.Item(263, 4).Value = 12909089
$ws.Cells.Item(264, 3).Value = 40048
$ws.Cells.Item(264, 4).Value = 56705812
$ws.Cells.Item(265, 3).Value = 4950
$ws.Cells.Item(265, 4).Value = 6384821
$ws.Cells.Item(266, 3).Value = 2662
$ws.Cells.Item(266, 4).Value = 3730046
$ws.Cells.Item(267, 3).Value = 10024
$ws.Cells.Item(267, 4).Value = 13271704
$ws.Cells.Item(268, 3).Value = 28922
$ws.Cells.Item(268, 4).Value = 38886657
$ws.Cells.Item(269, 3).Value = 14366
$ws.Cells.Item(269, 4).Value = 17956167
$ws.Cells.Item(271, 3).Value = 15822
$ws.Cells.Item(271, 4).Value = 18568856
$ws.Cells.Item(272, 3).Value = 32990
$ws.Cells.Item(272, 4).Value = 43690208
$ws.Cells.Item(273, 3).Value = 15488
$ws.Cells.Item(273, 4).Value = 18637049
$ws.Cells.Item(274, 3).Value = 39012
$ws.Cells.Item(274, 4).Value = 46800835
$ws.Cells.Item(275, 3).Value = 11338
$ws.Cells.Item(275, 4).Value = 15250352
$ws.Cells.Item(277, 3).Value = 17075
$ws.Cells.Item(277, 4).Value = 21195868
$ws.Cells.Item(278, 3).Value = 95
$ws.Cells.Item(278, 4).Value = 135212
$ws.Cells.Item(279, 3).Value = 650
$ws.Cells.Item(279, 4).Value = 811428
$ws.Cells.Item(280, 3).Value = 53729
$ws.Cells.Item(280, 4).Value = 76177656
$ws.Cells.Item(281, 3).Value = 66331
$ws.Cells.Item(281, 4).Value = 87172986
$ws.Cells.Item(282, 3).Value = 10693
$ws.Cells.Item(282, 4).Value = 14386775
$ws.Cells.Item(283, 3).Value = 47606
$ws.Cells.Item(283, 4).Value = 67547527
$ws.Cells.Item(284, 3).Value = 5934
$ws.Cells.Item(284, 4).Value = 7576264
$ws.Cells.Item(285, 3).Value = 2778
$ws.Cells.Item(285, 4).Value = 3842663
$ws.Cells.Item(286, 3).Value = 11410
$ws.Cells.Item(286, 4).Value = 15243674
$ws.Cells.Item(287, 3).Value = 32988
$ws.Cells.Item(287, 4).Value = 44724595
$ws.Cells.Item(288, 3).Value = 15799
$ws.Cells.Item(288, 4).Value = 19839948
$ws.Cells.Item(290, 3).Value = 19253
$ws.Cells.Item(290, 4).Value = 22637875
$ws.Cells.Item(291, 3).Value = 35834
$ws.Cells.Item(291, 4).Value = 47613898
$ws.Cells.Item(292, 3).Value = 17747
$ws.Cells.Item(292, 4).Value = 21285809
$ws.Cells.Item(293, 3).Value = 39861
$ws.Cells.Item(293, 4).Value = 46481799
$ws.Cells.Item(294, 3).Value = 3665
$ws.Cells.Item(294, 4).Value = 5125903
$ws.Cells.Item(296, 3).Value = 7038
$ws.Cells.Item(296, 4).Value = 8745815
$ws.Cells.Item(299, 3).Value = 17148
$ws.Cells.Item(299, 4).Value = 24591427
$ws.Cells.Item(300, 3).Value = 26188
$ws.Cells.Item(300, 4).Value = 34375076
$ws.Cells.Item(301, 3).Value = 4717
$ws.Cells.Item(301, 4).Value = 6603286
$ws.Cells.Item(302, 3).Value = 19464
$ws.Cells.Item(302, 4).Value = 27969166
$ws.Cells.Item(303, 3).Value = 2718
$ws.Cells.Item(303, 4).Value = 3483725
$ws.Cells.Item(304, 3).Value = 1639
$ws.Cells.Item(304, 4).Value = 2316041
$ws.Cells.Item(305, 3).Value = 5071
$ws.Cells.Item(305, 4).Value = 6793556
$ws.Cells.Item(306, 3).Value = 15188
$ws.Cells.Item(306, 4).Value = 20551790
$ws.Cells.Item(307, 3).Value = 5613
$ws.Cells.Item(307, 4).Value = 7227426
$ws.Cells.Item(308, 3).Value = 7187
$ws.Cells.Item(308, 4).Value = 8688926
$ws.Cells.Item(309, 3).Value = 18477
$ws.Cells.Item(309, 4).Value = 24436256
$ws.Cells.Item(310, 3).Value = 8543
$ws.Cells.Item(310, 4).Value = 10460473
$ws.Cells.Item(311, 3).Value = 19751
$ws.Cells.Item(311, 4).Value = 24354558
$ws.Cells.Item(312, 3).Value = 4905
$ws.Cells.Item(312, 4).Value = 6675789
$ws.Cells.Item(314, 3).Value = 16253
$ws.Cells.Item(314, 4).Value = 21107041
$ws.Cells.Item(317, 3).Value = 47157
$ws.Cells.Item(317, 4).Value = 67276730
$ws.Cells.Item(318, 3).Value = 70636
$ws.Cells.Item(318, 4).Value = 94116929
$ws.Cells.Item(319, 3).Value = 18083
$ws.Cells.Item(319, 4).Value = 25336184
$ws.Cells.Item(320, 3).Value = 45728
$ws.Cells.Item(320, 4).Value = 65994538
$ws.Cells.Item(321, 3).Value = 6443
$ws.Cells.Item(321, 4).Value = 8424438
$ws.Cells.Item(322, 3).Value = 3093
$ws.Cells.Item(322, 4).Value = 4362879
$ws.Cells.Item(323, 3).Value = 14076
$ws.Cells.Item(323, 4).Value = 19402010
$ws.Cells.Item(324, 3).Value = 35751
$ws.Cells.Item(324, 4).Value = 48656354
$ws.Cells.Item(325, 3).Value = 19062
$ws.Cells.Item(325, 4).Value = 24638215
$ws.Cells.Item(327, 3).Value = 20530
$ws.Cells.Item(327, 4).Value = 24973526
$ws.Cells.Item(328, 3).Value = 33015
$ws.Cells.Item(328, 4).Value = 44048429
$ws.Cells.Item(329, 3).Value = 16947
$ws.Cells.Item(329, 4).Value = 21034303
$ws.Cells.Item(330, 3).Value = 39537
$ws.Cells.Item(330, 4).Value = 47194964
